$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 12

$ws.Cells.Item($row, 1).Value = 13
$ws.Cells.Item($row, 2).Value = "'2026-02-16"
$ws.Cells.Item($row, 3).Value = "21:22:51"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "UP"
$ws.Cells.Item($row, 6).Value = 69433.7
$ws.Cells.Item($row, 7).Value = ""
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.75
$ws.Cells.Item($row, 12).Value = "Binance leading with 0.078% move"
$ws.Cells.Item($row, 13).Value = ""
$ws.Cells.Item($row, 14).Value = 0
